$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date Last Updated for Task 3 row moved forward
$ws.Range("D4").Value = 45414

# Task 5 & Task 6 rows gained a "Date Last Updated" value (copy date format from D5 first)
$ws.Range("D5").Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("D6").Value = 45445

$ws.Range("D5").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("D7").Value = 45445

# Task 7 status changed from "Pending" to "In Development" (copy style from G3 which already uses the "In Development" font)
$ws.Range("G3").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("B8").Value = "In Development"

# Task 11 gained a descriptive suffix
$ws.Range("A11").Value = "Task 11: Login Backend"

# Task 8 renamed from "Login Backend" to "Create Factory for dummy data"
$ws.Range("A9").Value = "Task 8: Create Factory for dummy data"

# Task 3 renamed from "Register Page" to "Register Modal"
$ws.Range("A4").Value = "Task 3: Register Modal"

# Update the active selection to match the saved view state
$ws.Range("C11").Select()
